$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 17:34"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1325791
$ws.Range("C4").Value = 4006
$ws.Range("D4").Value = 223939
$ws.Range("E4").Value = 1023090

# Row 7 - Reino Unido
$ws.Range("B7").Value = 215260
$ws.Range("C7").Value = 3896
$ws.Range("E7").Value = 183329
$ws.Range("G7").Value = 346
$ws.Range("H7").Value = 31587

# Row 10 - Alemania
$ws.Range("B10").Value = 170876
$ws.Range("C10").Value = 288
$ws.Range("E10").Value = 20066

# Row 17 - India
$ws.Range("B17").Value = 61356
$ws.Range("C17").Value = 1661
$ws.Range("D17").Value = 18672
$ws.Range("E17").Value = 40643
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = 2041

# Row 29 - Singapur
$ws.Range("D29").Value = 2296
$ws.Range("E29").Value = 20144
$ws.Range("F29").Value = 23

# Row 58 - Argelia
$ws.Range("B58").Value = 5558
$ws.Range("C58").Value = 189
$ws.Range("D58").Value = 2546
$ws.Range("E58").Value = 2518
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 494

# Row 70 - Grecia
$ws.Range("B70").Value = 2710
$ws.Range("C70").Value = 19
$ws.Range("E70").Value = 1185
$ws.Range("F70").Value = 28
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 151

# Row 132 - Montenegro
$ws.Range("D132").Value = 274
$ws.Range("E132").Value = 42
